$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Collapse duplicate "Rum" rows (rows 13-18 -> single row) ---
# Keep the first Rum row, retarget its key to the merged key name, then
# remove the other five duplicate rows.
$ws.Range("A13").Value = "ATA_ITEM_RUM"
$ws.Rows("14:18").Delete()

# --- Collapse duplicate "Orange" rows (rows 4,5,6 -> single row 4) ---
# Keep row 4, retarget its key to the merged key name, then remove rows 5 and 6.
$ws.Range("A4").Value = "ATA_ITEM_ORANGE"
$ws.Rows("5:6").Delete()

# --- Add the new "Tentacle Bow" entry as the new final row (row 14) ---
$ws.Range("A14").Value = "ATA_WEAPON_BOW_OF_ABYSS"
$ws.Range("C14").Value = "触手弓"
$ws.Range("B14").Value = "Tentacle Bow"

# --- Update selection to match the final authored state ---
$ws.Range("C14").Select()
